$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 121, pushing the existing rows 121:241 down to 122:242.
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new record (same shape as its
# neighbours: a weekly Perejil price observation for Vega Central Mapocho de
# Santiago, Región Metropolitana).
$ws.Cells.Item(121, 1).Value = 9
$ws.Cells.Item(121, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(121, 3).Value = "Metropolitana"
$ws.Cells.Item(121, 4).Value = 44512
$ws.Cells.Item(121, 5).Value = 13
$ws.Cells.Item(121, 6).Value = 100112044
$ws.Cells.Item(121, 7).Value = "Perejil"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 106
$ws.Cells.Item(121, 11).Value = 13000
$ws.Cells.Item(121, 12).Value = 14000
$ws.Cells.Item(121, 13).Value = 13500
$ws.Cells.Item(121, 14).Value = "`$/docena de atados"
$ws.Cells.Item(121, 15).Value = "Región Metropolitana"
$ws.Cells.Item(121, 16).Value = 4500
$ws.Cells.Item(121, 17).Value = 3
$ws.Cells.Item(121, 18).Value = "Hortaliza"
